$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 0
$ws.Range("B3").Value = 0.5
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.5

$ws.Range("B4").Select()
